$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 890 (pushes the existing rows 890:931 down to 891:932),
# growing the used range from A1:D931 to A1:D932.
$ws.Rows("890:890").Insert()

# Column A holds the date as literal text (not an Excel date serial), so force
# text format before assigning, then reset the cell style back to Normal so it
# matches the unstyled data cells around it.
$ws.Range("A890").NumberFormat = "@"
$ws.Range("A890").Value = "2026/02/28"
$ws.Range("A890").Style = "Normal"

$ws.Range("B890").Value = "土"
$ws.Range("C890").Value = 5
$ws.Range("D890").Value = 201
